# Updates cryptos list (price / 1h volume columns, plus a handful of
# row re-orderings where the coin at a given rank changed) per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Numeric-looking text values are written with a leading apostrophe so
# Excel keeps them as text (e.g. "1.00", "42.00") instead of coercing
# them to numbers and dropping the trailing zeros / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.742.09"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "3.461.65"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'592.06"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").Value = "'174.84"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("D9").Value = "'0.129"
$ws.Range("E9").Value = "  -3.57%  "

$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").Value = "'0.425"
$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("D12").Value = "4.050.03"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'30.61"
$ws.Range("E13").Value = "  +6.48%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "66.700.04"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -3.54%  "

$ws.Range("D17").Value = "3.427.98"
$ws.Range("E17").Value = "  -2.89%  "

$ws.Range("D18").Value = "'6.21"
$ws.Range("E18").Value = "  -2.39%  "

$ws.Range("D19").Value = "'14.24"
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "'386.36"
$ws.Range("E20").Value = "  -2.59%  "

$ws.Range("D21").Value = "'7.84"
$ws.Range("E21").Value = "  -1.84%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'72.57"
$ws.Range("E23").Value = "  -1.45%  "

$ws.Range("D24").Value = "'5.71"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'0.531"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  -2.95%  "

$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -2.62%  "

$ws.Range("D29").Value = "'0.994"
$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("D30").Value = "'6.06"
$ws.Range("E30").Value = "  -3.85%  "

$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -4.05%  "

$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("D33").Value = "'23.30"
$ws.Range("E33").Value = "  -3.21%  "

$ws.Range("D34").Value = "'7.21"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  -2.37%  "

$ws.Range("D36").Value = "'162.14"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("D37").Value = "'0.866"
$ws.Range("E37").Value = "  -3.74%  "

$ws.Range("D38").Value = "'1.91"
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("D39").Value = "'6.89"
$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'27.04"
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.59"
$ws.Range("E41").Value = "  -3.19%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'26.08"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.774.79"
$ws.Range("E43").Value = "  -1.62%  "

$ws.Range("D44").Value = "'0.0718"
$ws.Range("E44").Value = "  -3.88%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'42.00"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'339.32"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0298"
$ws.Range("E48").Value = "  -4.24%  "

$ws.Range("D49").Value = "'1.06"
$ws.Range("E49").Value = "  -3.41%  "

$ws.Range("D50").Value = "'32.97"
$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("D51").Value = "'0.103"
$ws.Range("E51").Value = "  -1.84%  "
